# "linking add_object to data_models"
# The sheet's single column was a list of training "Event Name" entries;
# it's being repurposed/renamed as a generic "data" column (header cell A1,
# which also drives the Table1 column header since A1 is inside Table1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "data"

# Reflect the author's last on-screen position: scrolled down so row 61 is
# at the top of the viewport, with H72 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H72").Select()
